$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row within the table at row 36 (pushes existing rows down)
$ws.Rows.Item(36).Insert()
$lo.Resize($ws.Range("A1:D50"))

# Populate the new row with the Paul Peak Trail data
$ws.Range("A36").Value = "Paul Peak Trail"
$ws.Range("B36").Value = 7.2
$ws.Range("C36").Value = 1630
$ws.Range("D36").Value = "moderate"
